$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (28 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2115646
$ws.Range("I132").Value = 2457929
$ws.Range("J132").Value = 4900
$ws.Range("K132").Value = 7373787
$ws.Range("L132").Value = 14700
$ws.Range("M132").Value = -7371257
$ws.Range("N132").Value = -19760
$ws.Range("H137").Value = 1481.6786
$ws.Range("I137").Value = 1246.5294
$ws.Range("J137").Value = 1845.091
$ws.Range("K137").Value = 3739.5882
$ws.Range("L137").Value = 5535.272999999999
$ws.Range("M137").Value = -1189.5882
$ws.Range("N137").Value = -10635.273
$ws.Range("H138").Value = 3636.2942
$ws.Range("I138").Value = 804.68085
$ws.Range("J138").Value = 7138.5527
$ws.Range("K138").Value = 2414.04255
$ws.Range("L138").Value = 21415.6581
$ws.Range("M138").Value = 2725.95745
$ws.Range("N138").Value = -31695.6581
$ws.Range("H141").Value = 2462.3977
$ws.Range("I141").Value = 1065.7903
$ws.Range("J141").Value = 6585.7144
$ws.Range("K141").Value = 3197.3709
$ws.Range("L141").Value = 19757.1432
$ws.Range("M141").Value = 1982.6291
$ws.Range("N141").Value = -30117.1432

# --- Sheet: ARM (32 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7257.826
$ws.Range("I32").Value = 3569.2307
$ws.Range("K32").Value = 3569.2307
$ws.Range("M32").Value = -3282.2307
$ws.Range("H74").Value = 859
$ws.Range("I74").Value = 783.62067
$ws.Range("J74").Value = 1171.2858
$ws.Range("K74").Value = 783.62067
$ws.Range("L74").Value = 1171.2858
$ws.Range("M74").Value = 90.37932999999998
$ws.Range("N74").Value = -2919.2858
$ws.Range("H77").Value = 859
$ws.Range("I77").Value = 783.62067
$ws.Range("J77").Value = 1171.2858
$ws.Range("K77").Value = 3918.10335
$ws.Range("L77").Value = 5856.429
$ws.Range("M77").Value = 449.8966499999997
$ws.Range("N77").Value = -14592.429
$ws.Range("H102").Value = 1365
$ws.Range("I102").Value = 1328.8235
$ws.Range("J102").Value = 1433.3334
$ws.Range("K102").Value = 1328.8235
$ws.Range("L102").Value = 1433.3334
$ws.Range("M102").Value = 293.1765
$ws.Range("N102").Value = -4677.3334
$ws.Range("H132").Value = 1697.326
$ws.Range("I132").Value = 1216.4546
$ws.Range("J132").Value = 2918
$ws.Range("K132").Value = 3649.3638
$ws.Range("L132").Value = 8754
$ws.Range("M132").Value = -1119.3638
$ws.Range("N132").Value = -13814

# --- Sheet: CRP (70 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3962.4546
$ws.Range("I16").Value = 2176.3333
$ws.Range("J16").Value = 12000
$ws.Range("K16").Value = 2176.3333
$ws.Range("L16").Value = 12000
$ws.Range("M16").Value = -1889.3333
$ws.Range("N16").Value = -12574
$ws.Range("H31").Value = 2631.5134
$ws.Range("I31").Value = 1370.08
$ws.Range("J31").Value = 5259.5
$ws.Range("K31").Value = 1370.08
$ws.Range("L31").Value = 5259.5
$ws.Range("M31").Value = -1075.08
$ws.Range("N31").Value = -5849.5
$ws.Range("H34").Value = 2631.5134
$ws.Range("I34").Value = 1370.08
$ws.Range("J34").Value = 5259.5
$ws.Range("K34").Value = 1370.08
$ws.Range("L34").Value = 5259.5
$ws.Range("M34").Value = -1168.08
$ws.Range("N34").Value = -5663.5
$ws.Range("H58").Value = 853.90625
$ws.Range("I58").Value = 660.4400000000001
$ws.Range("J58").Value = 1544.8572
$ws.Range("K58").Value = 660.4400000000001
$ws.Range("L58").Value = 1544.8572
$ws.Range("M58").Value = -457.4400000000001
$ws.Range("N58").Value = -1950.8572
$ws.Range("H99").Value = 2081.9167
$ws.Range("I99").Value = 1886.2222
$ws.Range("J99").Value = 2669
$ws.Range("K99").Value = 1886.2222
$ws.Range("L99").Value = 2669
$ws.Range("M99").Value = -388.2221999999999
$ws.Range("N99").Value = -5665
$ws.Range("H113").Value = 3962.4546
$ws.Range("I113").Value = 2176.3333
$ws.Range("J113").Value = 12000
$ws.Range("K113").Value = 2176.3333
$ws.Range("L113").Value = 12000
$ws.Range("M113").Value = -6.333299999999781
$ws.Range("N113").Value = -16340
$ws.Range("H126").Value = 2081.9167
$ws.Range("I126").Value = 1886.2222
$ws.Range("J126").Value = 2669
$ws.Range("K126").Value = 5658.6666
$ws.Range("L126").Value = 8007
$ws.Range("M126").Value = -3188.6666
$ws.Range("N126").Value = -12947
$ws.Range("H132").Value = 1565.5518
$ws.Range("I132").Value = 1056.7727
$ws.Range("J132").Value = 3164.5715
$ws.Range("K132").Value = 3170.3181
$ws.Range("L132").Value = 9493.7145
$ws.Range("M132").Value = -640.3181
$ws.Range("N132").Value = -14553.7145
$ws.Range("H134").Value = 1243.6744
$ws.Range("I134").Value = 1061.8
$ws.Range("J134").Value = 2039.375
$ws.Range("K134").Value = 3185.4
$ws.Range("L134").Value = 6118.125
$ws.Range("M134").Value = -650.3999999999996
$ws.Range("N134").Value = -11188.125
$ws.Range("H136").Value = 853.90625
$ws.Range("I136").Value = 660.4400000000001
$ws.Range("J136").Value = 1544.8572
$ws.Range("K136").Value = 1981.32
$ws.Range("L136").Value = 4634.571599999999
$ws.Range("M136").Value = 568.6799999999998
$ws.Range("N136").Value = -9734.571599999999

# --- Sheet: CUL (60 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1130.5
$ws.Range("J75").Value = 1130.5
$ws.Range("L75").Value = 3391.5
$ws.Range("N75").Value = -5387.5
$ws.Range("H78").Value = 1130.5
$ws.Range("J78").Value = 1130.5
$ws.Range("L78").Value = 10174.5
$ws.Range("N78").Value = -20158.5
$ws.Range("H81").Value = 1175
$ws.Range("J81").Value = 2050
$ws.Range("L81").Value = 6150
$ws.Range("N81").Value = -8396
$ws.Range("H82").Value = 4745
$ws.Range("I82").Value = 1490
$ws.Range("K82").Value = 4470
$ws.Range("M82").Value = -4064
$ws.Range("H84").Value = 1175
$ws.Range("J84").Value = 2050
$ws.Range("L84").Value = 18450
$ws.Range("N84").Value = -29682
$ws.Range("H85").Value = 4745
$ws.Range("I85").Value = 1490
$ws.Range("K85").Value = 4470
$ws.Range("M85").Value = -3066
$ws.Range("H86").Value = 625
$ws.Range("I86").Value = 700
$ws.Range("J86").Value = 600
$ws.Range("K86").Value = 2100
$ws.Range("L86").Value = 1800
$ws.Range("M86").Value = -914
$ws.Range("N86").Value = -4172
$ws.Range("H87").Value = 28493.857
$ws.Range("I87").Value = 1607
$ws.Range("J87").Value = 32975
$ws.Range("K87").Value = 4821
$ws.Range("L87").Value = 98925
$ws.Range("M87").Value = -3573
$ws.Range("N87").Value = -101421
$ws.Range("H88").Value = 3980
$ws.Range("J88").Value = 3980
$ws.Range("L88").Value = 11940
$ws.Range("N88").Value = -12796
$ws.Range("H89").Value = 625
$ws.Range("I89").Value = 700
$ws.Range("J89").Value = 600
$ws.Range("K89").Value = 6300
$ws.Range("L89").Value = 5400
$ws.Range("M89").Value = -372
$ws.Range("N89").Value = -17256
$ws.Range("H90").Value = 28493.857
$ws.Range("I90").Value = 1607
$ws.Range("J90").Value = 32975
$ws.Range("K90").Value = 14463
$ws.Range("L90").Value = 296775
$ws.Range("M90").Value = -8223
$ws.Range("N90").Value = -309255
$ws.Range("H91").Value = 3980
$ws.Range("J91").Value = 3980
$ws.Range("L91").Value = 11940
$ws.Range("N91").Value = -14904

# --- Sheet: LTW (35 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2545.3076
$ws.Range("I7").Value = 2958.8
$ws.Range("J7").Value = 2286.875
$ws.Range("K7").Value = 2958.8
$ws.Range("L7").Value = 2286.875
$ws.Range("M7").Value = -2846.8
$ws.Range("N7").Value = -2510.875
$ws.Range("H40").Value = 2133.8147
$ws.Range("I40").Value = 1910.45
$ws.Range("J40").Value = 2772
$ws.Range("K40").Value = 1910.45
$ws.Range("L40").Value = 2772
$ws.Range("M40").Value = -1774.45
$ws.Range("N40").Value = -3044
$ws.Range("H93").Value = 4817.4
$ws.Range("I93").Value = 7001
$ws.Range("J93").Value = 1961.9231
$ws.Range("K93").Value = 7001
$ws.Range("L93").Value = 1961.9231
$ws.Range("M93").Value = -5753
$ws.Range("N93").Value = -4457.9231
$ws.Range("H126").Value = 2545.3076
$ws.Range("I126").Value = 2958.8
$ws.Range("J126").Value = 2286.875
$ws.Range("K126").Value = 8876.400000000001
$ws.Range("L126").Value = 6860.625
$ws.Range("M126").Value = -6406.400000000001
$ws.Range("N126").Value = -11800.625
$ws.Range("H132").Value = 1337.6562
$ws.Range("I132").Value = 920.4286
$ws.Range("J132").Value = 2700.6
$ws.Range("K132").Value = 2761.2858
$ws.Range("L132").Value = 8101.799999999999
$ws.Range("M132").Value = -231.2857999999997
$ws.Range("N132").Value = -13161.8

# --- Sheet: WVR (21 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1778.7727
$ws.Range("I122").Value = 1512.3
$ws.Range("J122").Value = 2000.8334
$ws.Range("K122").Value = 4536.9
$ws.Range("L122").Value = 6002.5002
$ws.Range("M122").Value = -2086.9
$ws.Range("N122").Value = -10902.5002
$ws.Range("H132").Value = 1002.0238
$ws.Range("I132").Value = 692.4483
$ws.Range("J132").Value = 1692.6154
$ws.Range("K132").Value = 2077.3449
$ws.Range("L132").Value = 5077.8462
$ws.Range("M132").Value = 452.6550999999999
$ws.Range("N132").Value = -10137.8462
$ws.Range("H136").Value = 2820.4736
$ws.Range("I136").Value = 841.619
$ws.Range("J136").Value = 8361.267
$ws.Range("K136").Value = 2524.857
$ws.Range("L136").Value = 25083.801
$ws.Range("M136").Value = 25.14300000000003
$ws.Range("N136").Value = -30183.801

Write-Output "Applied all cell updates"